# Update rh_license_approvals field mapping: rename columns to match the
# new ER (replacement_id -> solicitation_id, next_approved -> next_approved_user_id)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B12").Value = "solicitation_id"
$ws.Range("B13").Value = "next_approved_user_id"
$ws.Range("B15").Value = "solicitation_id"

# Keep selection consistent with the authored workbook state
$ws.Range("B19").Select()
